# TC09_C3DC_phs002431_TrtmntType-Imunothrpy.xlsx
# "Updated C3DC phs002431 testcases"
#
# The TreatmentTab query (Sheet1!B5) needs its WHERE clause tightened so it
# also excludes rows with a NULL treatment id, matching the other queries'
# refreshed filters on this sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$cell = $ws.Range("B5")
$oldText = $cell.Value()

$oldWhere = "std.dbgap_accession = 'phs002431' AND trt.treatment_type = 'Immunotherapy'"
$newWhere = "std.dbgap_accession = 'phs002431' AND trt.treatment_type = 'Immunotherapy' AND trt.treatment_id IS NOT NULL"

$cell.Value = $oldText.Replace($oldWhere, $newWhere)

# Re-apply the cell's formatting (12pt Calibri, wrapped text) so the edited
# cell picks up a freshly written style, same as happens when Excel re-saves
# a retyped/pasted cell.
$cell.Font.Name = "Calibri"
$cell.Font.Size = 12
$cell.WrapText = $true

# Leave the selection/viewport on the row that was just edited.
[void]$ws.Range("C5").Select()

Write-Output "Updated Sheet1!B5 treatment query WHERE clause"
